# Adds a new "Sizes" column (column N) to every product sheet in the
# workbook, populated with size information for each row.
$wb = $excel.ActiveWorkbook

# --- Bags_df ---
$ws = $wb.Worksheets.Item("Bags_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "N/A"
$ws.Range("N3").Value = "N/A"
$ws.Range("N4").Value = "N/A"
$ws.Range("N5").Value = "N/A"
$ws.Range("N6").Value = "N/A"
$ws.Range("N7").Value = "N/A"

# --- Skate_df ---
$ws = $wb.Worksheets.Item("Skate_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "129 ,139 ,149 "
$ws.Range("N3").Value = "54MM"
$ws.Range("N4").Value = "8 1/2"
$ws.Range("N5").Value = "8 1/8"

# --- Shirts_df ---
$ws = $wb.Worksheets.Item("Shirts_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "Small,Medium,Large,XLarge"
$ws.Range("N3").Value = "Small,Medium,Large,XLarge"
$ws.Range("N4").Value = "Small,Medium,Large,XLarge"
$ws.Range("N5").Value = "Small,Medium,Large,XLarge"

# --- Pants_df ---
$ws = $wb.Worksheets.Item("Pants_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "Small,Medium,Large,XLarge"
$ws.Range("N3").Value = "30 ,32 ,34 ,36 ,38 "
$ws.Range("N4").Value = "30 ,32 ,34 ,36 ,38 "
$ws.Range("N5").Value = "30 ,32 ,34 ,36 ,38 "
$ws.Range("N6").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N7").Value = "30 ,32 ,34 ,36 ,38 "
$ws.Range("N8").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N9").Value = "30 ,32 ,34 ,36 ,38 "
$ws.Range("N10").Value = "30 ,32 ,34 ,36 ,38 "
$ws.Range("N11").Value = "30 ,32 ,34 ,36 ,38 "
$ws.Range("N12").Value = "30 ,32 ,34 ,36 ,38 "

# --- Shorts_df ---
$ws = $wb.Worksheets.Item("Shorts_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "Small,Medium,Large,XLarge"
$ws.Range("N3").Value = "Small,Medium,Large,XLarge,XXL"

# --- Tops_Sweaters_df ---
$ws = $wb.Worksheets.Item("Tops_Sweaters_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "Small,Medium,Large,XLarge"
$ws.Range("N3").Value = "Small,Medium,Large,XLarge"
$ws.Range("N4").Value = "Small,Medium,Large,XLarge"
$ws.Range("N5").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N6").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N7").Value = "Small,Medium,Large,XLarge"
$ws.Range("N8").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N9").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N10").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N11").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N12").Value = "Small,Medium,Large,XLarge,XXL"

# --- T_Shirts_df ---
$ws = $wb.Worksheets.Item("T_Shirts_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N3").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N4").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N5").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N6").Value = "Small,Medium,Large,XLarge,XXL"

# --- Jackets_df ---
$ws = $wb.Worksheets.Item("Jackets_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "Small,Medium,Large,XLarge"
$ws.Range("N3").Value = "Small,Medium,Large,XLarge"
$ws.Range("N4").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N5").Value = "Small,Medium,Large,XLarge"
$ws.Range("N6").Value = "Small,Medium,Large,XLarge"
$ws.Range("N7").Value = "Small,Medium,Large,XLarge,XXL"

# --- Sweatshirts_df ---
$ws = $wb.Worksheets.Item("Sweatshirts_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N3").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N4").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N5").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N6").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N7").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N8").Value = "Small,Medium,Large,XLarge,XXL"

# --- Hats_df ---
$ws = $wb.Worksheets.Item("Hats_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "N/A"
$ws.Range("N3").Value = "N/A"
$ws.Range("N4").Value = "N/A"
$ws.Range("N5").Value = "N/A"
$ws.Range("N6").Value = "S/M,M/L"
$ws.Range("N7").Value = "N/A"
$ws.Range("N8").Value = "N/A"

# --- Accessories_df ---
$ws = $wb.Worksheets.Item("Accessories_df")
$ws.Range("N1").Value = "Sizes"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N2").Value = "N/A"
$ws.Range("N3").Value = "N/A"
$ws.Range("N4").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N5").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N6").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N7").Value = "Small,Medium,Large,XLarge,XXL"
$ws.Range("N8").Value = "N/A"

$excel.CutCopyMode = $false
